$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header text updates (shared-string rich-text runs) ----
$ws.Range("A8").Characters(21, 2).Text = "16"
$ws.Range("C9").Characters(27, 8).Text = "4/14/2025"
$ws.Range("C9").Characters(47, 9).Text = "4/20/2025"

# ---- Cells changing from numeric style to shared-string "N/A" style (style 13) ----
# Use D14 (style 13, text "0") and E14 (style 13, text "***.*") as format+value donors.
$ws.Range("D14").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("C14").PasteSpecial(-4163)
$ws.Range("D14").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("G15").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("H15").PasteSpecial(-4163)
$ws.Range("D14").Copy()
$ws.Range("G27").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("G27").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("H27").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("H27").PasteSpecial(-4163)
$ws.Range("D14").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("C29").PasteSpecial(-4163)
$ws.Range("D14").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("C30").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# ---- Cells changing from shared-string "N/A" style to numeric style (style 14 or 15) ----
# Use F14 (style 14) and L14 (style 15) as format donors, then set the numeric value.
$ws.Range("F14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D22").Value = 1
$ws.Range("L14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E22").Value = -100
$ws.Range("F14").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("D23").Value = 1
$ws.Range("L14").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("E23").Value = -100
$ws.Range("F14").Copy()
$ws.Range("G23").PasteSpecial(-4122)
$ws.Range("G23").Value = 1
$ws.Range("L14").Copy()
$ws.Range("H23").PasteSpecial(-4122)
$ws.Range("H23").Value = 100
$ws.Range("F14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C28").Value = 1
$ws.Range("F14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D28").Value = 2
$ws.Range("L14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E28").Value = -50
$excel.CutCopyMode = $false

# ---- Plain numeric value updates (style unchanged) ----
$ws.Range("M15").Value = -50
$ws.Range("N15").Value = -82.608695652173
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 22
$ws.Range("H16").Value = -50
$ws.Range("I16").Value = 52
$ws.Range("J16").Value = 111
$ws.Range("K16").Value = -53.153153153153
$ws.Range("L16").Value = -42.222222222222
$ws.Range("M16").Value = -56.666666666666
$ws.Range("N16").Value = -87.878787878787
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 10
$ws.Range("F17").Value = 24
$ws.Range("G17").Value = 38
$ws.Range("H17").Value = -36.842105263157
$ws.Range("I17").Value = 93
$ws.Range("J17").Value = 128
$ws.Range("K17").Value = -27.34375
$ws.Range("L17").Value = -19.130434782608
$ws.Range("M17").Value = 9.411764705882
$ws.Range("N17").Value = -64.367816091954
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 50
$ws.Range("I18").Value = 44
$ws.Range("J18").Value = 73
$ws.Range("K18").Value = -39.726027397260
$ws.Range("L18").Value = -37.142857142857
$ws.Range("M18").Value = -50.561797752809
$ws.Range("N18").Value = -89.189189189189
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -8.333333333333
$ws.Range("F19").Value = 29
$ws.Range("G19").Value = 57
$ws.Range("H19").Value = -49.122807017543
$ws.Range("I19").Value = 144
$ws.Range("J19").Value = 185
$ws.Range("K19").Value = -22.162162162162
$ws.Range("L19").Value = -33.333333333333
$ws.Range("M19").Value = 61.797752808988
$ws.Range("N19").Value = -6.493506493506
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 400
$ws.Range("F20").Value = 12
$ws.Range("H20").Value = -20
$ws.Range("I20").Value = 43
$ws.Range("J20").Value = 58
$ws.Range("K20").Value = -25.862068965517
$ws.Range("L20").Value = -10.416666666666
$ws.Range("M20").Value = -4.444444444444
$ws.Range("N20").Value = -86.217948717948
$ws.Range("C21").Value = 28
$ws.Range("D21").Value = 29
$ws.Range("E21").Value = -3.448275862068
$ws.Range("F21").Value = 88
$ws.Range("G21").Value = 147
$ws.Range("H21").Value = -40.136054421768
$ws.Range("I21").Value = 384
$ws.Range("J21").Value = 564
$ws.Range("K21").Value = -31.914893617021
$ws.Range("L21").Value = -29.541284403669
$ws.Range("M21").Value = -12.328767123287
$ws.Range("N21").Value = -76
$ws.Range("G22").Value = 5
$ws.Range("J22").Value = 13
$ws.Range("K22").Value = -76.923076923076
$ws.Range("L22").Value = -72.727272727272
$ws.Range("F23").Value = 2
$ws.Range("J23").Value = 9
$ws.Range("K23").Value = 22.222222222222
$ws.Range("L23").Value = 10
$ws.Range("M23").Value = 120
$ws.Range("C24").Value = 25
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = 38.888888888888
$ws.Range("F24").Value = 81
$ws.Range("G24").Value = 56
$ws.Range("H24").Value = 44.642857142857
$ws.Range("I24").Value = 288
$ws.Range("J24").Value = 279
$ws.Range("K24").Value = 3.225806451612
$ws.Range("L24").Value = 7.462686567164
$ws.Range("M24").Value = 41.871921182266
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 50
$ws.Range("F25").Value = 27
$ws.Range("G25").Value = 13
$ws.Range("H25").Value = 107.692307692308
$ws.Range("I25").Value = 78
$ws.Range("J25").Value = 67
$ws.Range("K25").Value = 16.417910447761
$ws.Range("L25").Value = 62.5
$ws.Range("C26").Value = 15
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 114.285714285714
$ws.Range("F26").Value = 59
$ws.Range("G26").Value = 38
$ws.Range("H26").Value = 55.263157894736
$ws.Range("I26").Value = 213
$ws.Range("J26").Value = 188
$ws.Range("K26").Value = 13.297872340425
$ws.Range("L26").Value = 15.135135135135
$ws.Range("M26").Value = 0.471698113207
$ws.Range("L27").Value = -54.545454545454
$ws.Range("F28").Value = 2
$ws.Range("H28").Value = -33.333333333333
$ws.Range("I28").Value = 13
$ws.Range("J28").Value = 19
$ws.Range("K28").Value = -31.578947368421
$ws.Range("L28").Value = -13.333333333333
$ws.Range("F29").Value = 2
$ws.Range("M29").Value = -30
$ws.Range("N29").Value = -88.888888888888
$ws.Range("F30").Value = 2
$ws.Range("M30").Value = -33.333333333333
$ws.Range("N30").Value = -90.322580645161
